# Folder Restructure + Added Al Santillan Time Sheets.
#
# Adds a new "TA" meeting/time-sheet block (rows 12-18) to the
# "Spring 2021" sheet, mirroring the existing TEAM block in rows 3-9,
# and tidies up the header-row styling (row 3/4 empty-style cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spring 2021")

# Shrink/reposition the saved window (best-effort; cosmetic window chrome)
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 460
$win.Width = 25600
$win.Height = 15540

# --- Row 3: drop the stray C3 style-only cell (keep D3:M3 as-is) ---
$ws.Range("C3").Clear() | Out-Null

# --- Row 4: extend the centred empty-style band from C4 out to D4:M4 ---
$ws.Range("D4:M4").HorizontalAlignment = -4108   # xlCenter (matches style used by C3:M3/D3:M3)

# --- New time-sheet block mirroring rows 3-9, placed at rows 12-18 ---

# Date header (same style as B3)
$ws.Range("B12").Value = 44229
$ws.Range("B12").HorizontalAlignment = $ws.Range("B3").HorizontalAlignment
$ws.Range("B12").NumberFormat = $ws.Range("B3").NumberFormat

# "Meeting Type" / "TA" header row (mirrors A4/B4, but TA instead of TEAM)
$ws.Range("A13").Value = $ws.Range("A4").Text
$ws.Range("B13").Value = "TA"
$ws.Range("B13").HorizontalAlignment = $ws.Range("B4").HorizontalAlignment

# Attendee rows (mirrors A5:B9 -> A14:B18)
$ws.Range("A14").Value = $ws.Range("A5").Text
$ws.Range("B14").Interior.Color = $ws.Range("B5").Interior.Color

$ws.Range("A15").Value = $ws.Range("A6").Text
$ws.Range("B15").Interior.Color = $ws.Range("B6").Interior.Color

$ws.Range("A16").Value = $ws.Range("A7").Text
$ws.Range("B16").Interior.Color = $ws.Range("B7").Interior.Color

$ws.Range("A17").Value = $ws.Range("A8").Text
$ws.Range("B17").Interior.Color = $ws.Range("B8").Interior.Color

$ws.Range("A18").Value = $ws.Range("A9").Text
$ws.Range("B18").Interior.Color = $ws.Range("B9").Interior.Color

# Restore selection to where the author left off editing
$ws.Range("D13").Select() | Out-Null
